$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data (columns B, E-X, AA-AD) between row 11 and row 12,
# leaving columns A, C, D, I, K, Y, Z unchanged.

$cols = @("B","E","F","G","H","J","L","M","N","O","P","Q","R","S","T","U","V","W","X","AA","AB","AC","AD")

foreach ($col in $cols) {
    $addr11 = "$col`11"
    $addr12 = "$col`12"
    $v11 = $ws.Range($addr11).Value2
    $v12 = $ws.Range($addr12).Value2
    $ws.Range($addr11).Value2 = $v12
    $ws.Range($addr12).Value2 = $v11
}
